# Weekly NYPD CompStat refresh: new crime data collected for the
# reporting week of 1/22/2024 - 1/28/2024 (Volume 31, Number 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a number into a cell while copying the number-format /
# font / alignment from a "donor" cell that already carries the style
# we want the destination to end up with (reuses existing style table
# entries instead of inventing new ones).
# ---------------------------------------------------------------------
function Set-NumCell($ref, $value, $donorStyleRef) {
    $ws.Range($donorStyleRef).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $value
}

# ---------------------------------------------------------------------
# Helper: write a text label ("0" or "***.*") into a cell that
# currently holds a number, forcing it to be stored as text, then
# restore the donor's number format/alignment so the final style
# matches a normal "N/A" label cell.
# ---------------------------------------------------------------------
function Set-TextCell($ref, $text, $donorStyleRef) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($donorStyleRef).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Header: bump the report Volume/Number and the covered week's dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# ---------------------------------------------------------------------
# Murder row (14) - 14-year % chg column now has data instead of N/A.
# ---------------------------------------------------------------------
Set-NumCell "M14" 0 "E16"

# ---------------------------------------------------------------------
# Robbery row (16)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 16.666666666666
$ws.Range("I16").Value = 21
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = 16.666666666666
$ws.Range("L16").Value = 10.526315789473
$ws.Range("M16").Value = -4.545454545454
$ws.Range("N16").Value = -68.181818181818

# ---------------------------------------------------------------------
# Fel. Assault row (17)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 350
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 17.391304347826
$ws.Range("I17").Value = 27
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = 17.391304347826
$ws.Range("L17").Value = 22.727272727272
$ws.Range("M17").Value = 125
$ws.Range("N17").Value = -27.027027027027

# ---------------------------------------------------------------------
# Burglary row (18) - Week-to-Date 2024 count dropped to zero, so the
# 2024 column and its % chg become text "N/A" labels.
# ---------------------------------------------------------------------
Set-TextCell "C18" "0" "C14"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = -40
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = -40
$ws.Range("L18").Value = -75
$ws.Range("M18").Value = -75
$ws.Range("N18").Value = -92.307692307692

# ---------------------------------------------------------------------
# Gr. Larceny row (19)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 300
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = 12.5
$ws.Range("I19").Value = 27
$ws.Range("J19").Value = 24
$ws.Range("K19").Value = 12.5
$ws.Range("L19").Value = -18.181818181818
$ws.Range("M19").Value = 68.75
$ws.Range("N19").Value = 17.391304347826

# ---------------------------------------------------------------------
# G.L.A. row (20) - Week-to-Date now has real counts instead of N/A.
# ---------------------------------------------------------------------
Set-NumCell "C20" 1 "C16"
Set-NumCell "D20" 1 "C16"
Set-NumCell "E20" 0 "E16"
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 400
$ws.Range("I20").Value = 5
Set-NumCell "J20" 1 "C16"
Set-NumCell "K20" 400 "E16"
$ws.Range("L20").Value = 25
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -80.769230769230

# ---------------------------------------------------------------------
# TOTAL row (21)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 57.142857142857
$ws.Range("F21").Value = 85
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 18.055555555555
$ws.Range("I21").Value = 85
$ws.Range("J21").Value = 72
$ws.Range("K21").Value = 18.055555555555
$ws.Range("L21").Value = -6.593406593406
$ws.Range("M21").Value = 21.428571428571
$ws.Range("N21").Value = -55.729166666666

# ---------------------------------------------------------------------
# Transit row (22) - 28-Day 2023 count / % chg become N/A.
# ---------------------------------------------------------------------
$ws.Range("F22").Value = 4
Set-TextCell "G22" "0" "C14"
Set-TextCell "H22" "***.*" "C14"

# ---------------------------------------------------------------------
# Housing row (23) - Week-to-Date 2023 count / % chg become N/A.
# ---------------------------------------------------------------------
Set-TextCell "D23" "0" "C14"
Set-TextCell "E23" "***.*" "C14"
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 87.5
$ws.Range("I23").Value = 15
$ws.Range("K23").Value = 87.5
$ws.Range("L23").Value = 150
$ws.Range("M23").Value = 87.5

# ---------------------------------------------------------------------
# Petit Larceny row (24)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -60
$ws.Range("F24").Value = 56
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = -34.883720930232
$ws.Range("I24").Value = 56
$ws.Range("J24").Value = 86
$ws.Range("K24").Value = -34.883720930232
$ws.Range("L24").Value = 27.272727272727
$ws.Range("M24").Value = -30.864197530864

# ---------------------------------------------------------------------
# Misd. Assault row (25)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 71.428571428571
$ws.Range("F25").Value = 41
$ws.Range("H25").Value = 2.5
$ws.Range("I25").Value = 41
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = 2.5
$ws.Range("L25").Value = -2.380952380952
$ws.Range("M25").Value = 13.888888888888

# ---------------------------------------------------------------------
# Other Sex Crimes row (27) - Week-to-Date 2024 count now zero.
# ---------------------------------------------------------------------
Set-TextCell "C27" "0" "C14"
Set-NumCell "D27" 1 "C16"
Set-NumCell "E27" -100 "E16"
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 25

# ---------------------------------------------------------------------
# Shooting Vic. row (28)
# ---------------------------------------------------------------------
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = 3

# ---------------------------------------------------------------------
# Shooting Inc. row (29)
# ---------------------------------------------------------------------
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 3

# ---------------------------------------------------------------------
# Hate Crimes row (30) - 2-Year % chg now has data instead of N/A.
# ---------------------------------------------------------------------
Set-NumCell "L30" -100 "E16"
